$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# B7 (next to "Experimental") gets the text value "true".
# A plain Value assignment of the literal "true" would be auto-coerced to
# the Excel Boolean TRUE (same as typing true into a cell and having it
# autocomplete to a boolean). Going through a formula that evaluates to the
# string "true" and then converting that formula to its static result via
# Copy/PasteSpecial(xlPasteValues) keeps it as literal text.
$ws.Range("B7").Formula = "=""true"""
$ws.Range("B7").Copy()
$ws.Range("B7").PasteSpecial(-4163)

# B8 (next to "Date") gets the updated ISO-8601 timestamp.
$ws.Range("B8").Value = "2023-02-16T14:43:10-06:00"
